# Updated cryptos list on Wed Oct 18 02:35:50 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) columns for
# each coin row with newly scraped values. Both columns are stored as plain
# text in the sheet (e.g. "28.490.52", "  +0.38%  "), so for column D - whose
# new values often still *look* numeric (e.g. "212.17", "0.0590") - the COM
# Value setter is guarded with a temporary Text ("@") number format to stop
# Excel from silently coercing the string into a real number (which would
# drop significant trailing zeros). The cell style is then restored to
# "Normal" so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.490.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.571.03'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.47%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.17'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.25%  '
$ws.Range("E6").Value = '  -0.37%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '46.22'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.03'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.38%  '
$ws.Range("E10").Value = '  -1.86%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0590'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0887'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.793.98'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.566.82'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.520'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.470.86'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.30%  '
$ws.Range("E17").Value = '  -2.29%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.11'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '231.23'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0689'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.05%  '
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("E23").Value = '  -5.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.12'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.76'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.28%  '
$ws.Range("E28").Value = '  -2.57%  '
$ws.Range("E29").Value = '  -3.70%  '
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0479'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.56%  '
$ws.Range("E32").Value = '  -3.50%  '
$ws.Range("E33").Value = '  -1.45%  '
$ws.Range("E34").Value = '  -2.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.391.24'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.01%  '
$ws.Range("E37").Value = '  -3.78%  '
$ws.Range("E38").Value = '  +0.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.63'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.67%  '
$ws.Range("E40").Value = '  -0.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.520'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.89'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.63%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.787'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.41%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0467'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.90%  '
$ws.Range("E46").Value = '  -4.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.970'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '62.79'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.706.65'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.33'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.27%  '
$ws.Range("E51").Value = '  -0.96%  '
